$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): update B3 value, clear D3 value
$ws.Range("B3").Value = 21051850.29180706
$ws.Range("D3").Value = ""

# Row 4 (Methanol): update C4 value
$ws.Range("C4").Value = 8473.380451641917

# Row 5 (Ammonia): update C5 value
$ws.Range("C5").Value = 12713.77260147078

# Row 7: rename "Other" -> "Biogas", update D7 value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 9920.530001058689

# New row 8: "Other" entry (same formatting as A7), with D8 value
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 6768.185524167529
